$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "'2"
$ws.Range("D3").Value = "'23.11"
$ws.Range("G3").Value = "'2"
$ws.Range("G4").Value = "'2"
$ws.Range("D5").Value = "'0.05968"
$ws.Range("G5").Value = "'2"
$ws.Range("G6").Value = "'2"
$ws.Range("D7").Value = "'6.529"
$ws.Range("G7").Value = "'2"
$ws.Range("D8").Value = "'0.8129"
$ws.Range("G8").Value = "'2"
$ws.Range("D9").Value = "'0.9291"
$ws.Range("G9").Value = "'2"
$ws.Range("D10").Value = "'0.1429"
$ws.Range("G10").Value = "'2"
$ws.Range("D11").Value = "'0.07391"
$ws.Range("G11").Value = "'2"
$ws.Range("D12").Value = "'0.03280"
$ws.Range("G12").Value = "'2"
$ws.Range("D13").Value = "'0.03079"
$ws.Range("G13").Value = "'2"
$ws.Range("D14").Value = "'0.09354"
$ws.Range("G14").Value = "'2"
$ws.Range("D15").Value = "'3.849"
$ws.Range("G15").Value = "'2"
$ws.Range("D16").Value = "'0.001576"
$ws.Range("G16").Value = "'2"
$ws.Range("D17").Value = "'0.04702"
$ws.Range("G17").Value = "'2"
$ws.Range("D18").Value = "'0.0005915"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("G18").Value = "'2"
$ws.Range("D19").Value = "'0.005946"
$ws.Range("G19").Value = "'2"
$ws.Range("D20").Value = "'0.001270"
$ws.Range("G20").Value = "'2"
$ws.Range("D21").Value = "'0.004912"
$ws.Range("G21").Value = "'2"
$ws.Range("D22").Value = "'0.00006806"
$ws.Range("G22").Value = "'2"
$ws.Range("D23").Value = "'3.594"
$ws.Range("G23").Value = "'2"
$ws.Range("G24").Value = "'2"
$ws.Range("D25").Value = "'0.3235"
$ws.Range("G25").Value = "'2"
$ws.Range("D26").Value = "'0.1332"
$ws.Range("G26").Value = "'2"
$ws.Range("D27").Value = "'0.0003702"
$ws.Range("G27").Value = "'2"
$ws.Range("G28").Value = "'2"
$ws.Range("G29").Value = "'2"
$ws.Range("G30").Value = "'2"
$ws.Range("G31").Value = "'2"
$ws.Range("G32").Value = "'2"
$ws.Range("G33").Value = "'2"
$ws.Range("G34").Value = "'2"
$ws.Range("G35").Value = "'2"
$ws.Range("G36").Value = "'2"
$ws.Range("G37").Value = "'2"
$ws.Range("G38").Value = "'2"
$ws.Range("G39").Value = "'2"
$ws.Range("D40").Value = "'0.03963"
$ws.Range("G40").Value = "'2"
$ws.Range("D41").Value = "'0.006406"
$ws.Range("G41").Value = "'2"
$ws.Range("D42").Value = "'0.1078"
$ws.Range("G42").Value = "'2"
$ws.Range("G43").Value = "'2"
$ws.Range("D44").Value = "'0.01053"
$ws.Range("G44").Value = "'2"
$ws.Range("D45").Value = "'0.00005226"
$ws.Range("G45").Value = "'2"
$ws.Range("G46").Value = "'2"
$ws.Range("D47").Value = "'0.7256"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINWorstin24h"
$ws.Range("G47").Value = "'2"
$ws.Range("D48").Value = "'0.002327"
$ws.Range("G48").Value = "'2"
$ws.Range("G49").Value = "'2"
$ws.Range("G50").Value = "'2"
$ws.Range("G51").Value = "'2"
